# AutoJoin-style update: fix a couple of accented names and append a new
# "ADICIONADO" (added) record highlighted in green, mirroring how the sheet
# already highlights edited ("ALTERADO") rows in yellow/red.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (missing accents) -------------------------------
$ws.Range("F3").Value = "MARIO ANDRÉ"
$ws.Range("K3").Value = "JOÃO PESSOA"
$ws.Range("G6").Value = "JULIA MÊRCEDES"
$ws.Range("K6").Value = "JOÃO PESSOA"
$ws.Range("K15").Value = "JOÃO PESSOA"

# --- Append new row 16 --------------------------------------------------
# The sheet's existing rows store every value (even numeric-looking ones
# like MATRICULA/CPF/CEP/DDD/...) as plain text, so force Text format on
# the new row before writing to avoid Excel auto-converting to numbers.
$newRow = 16
$newRowRange = $ws.Range("A" + $newRow + ":Q" + $newRow)
$newRowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "15"
$ws.Cells.Item($newRow, 2).Value = "MARIA LUIZA"
$ws.Cells.Item($newRow, 3).Value = "SUPORTE TÉCNICO"
$ws.Cells.Item($newRow, 4).Value = "94585301234"
$ws.Cells.Item($newRow, 5).Value = "RUA J"
$ws.Cells.Item($newRow, 6).Value = "CARLOS ALBERTO"
$ws.Cells.Item($newRow, 7).Value = "JULIA FERNANDES"
$ws.Cells.Item($newRow, 8).Value = "566"
$ws.Cells.Item($newRow, 10).Value = "BAIRRO A"
$ws.Cells.Item($newRow, 11).Value = "JOÃO PESSOA"
$ws.Cells.Item($newRow, 12).Value = "PB"
$ws.Cells.Item($newRow, 13).Value = "52210901"
$ws.Cells.Item($newRow, 14).Value = "83"
$ws.Cells.Item($newRow, 15).Value = "935859334"
$ws.Cells.Item($newRow, 17).Value = "ADICIONADO"

# Restore General format now that the text values are locked in (matches
# the rest of the sheet, which stores numeric-looking values as text but
# keeps the General number format).
$newRowRange.NumberFormat = "General"

# Highlight the newly added row in green, same pattern used for the
# existing yellow/red "ALTERADO" rows.
$newRange = $ws.Range("A" + $newRow + ":Q" + $newRow)
$newRange.Interior.Color = 4709952
